$d = $word.ActiveDocument

$d.Content.Find.Execute("110÷6=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "421÷3=140, 1", 2)
$d.Content.Find.Execute("649÷5=129, 4", $true, $false, $false, $false, $false, $true, 1, $false, "267÷5=53, 2", 2)
$d.Content.Find.Execute("540÷5=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "643÷9=71, 4", 2)
$d.Content.Find.Execute("194÷7=27, 5", $true, $false, $false, $false, $false, $true, 1, $false, "330÷9=36, 6", 2)
$d.Content.Find.Execute("563÷5=112, 3", $true, $false, $false, $false, $false, $true, 1, $false, "190÷3=63, 1", 2)
$d.Content.Find.Execute("223÷7=31, 6", $true, $false, $false, $false, $false, $true, 1, $false, "324÷3=108, 0", 2)
$d.Content.Find.Execute("454÷2=227, 0", $true, $false, $false, $false, $false, $true, 1, $false, "607÷4=151, 3", 2)
$d.Content.Find.Execute("234÷9=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "869÷2=434, 1", 2)
$d.Content.Find.Execute("407÷9=45, 2", $true, $false, $false, $false, $false, $true, 1, $false, "648÷2=324, 0", 2)
$d.Content.Find.Execute("548÷7=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "623÷3=207, 2", 2)
$d.Content.Find.Execute("577÷9=64, 1", $true, $false, $false, $false, $false, $true, 1, $false, "946÷9=105, 1", 2)
$d.Content.Find.Execute("643÷2=321, 1", $true, $false, $false, $false, $false, $true, 1, $false, "845÷4=211, 1", 2)
$d.Content.Find.Execute("271÷6=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "711÷2=355, 1", 2)
$d.Content.Find.Execute("652÷6=108, 4", $true, $false, $false, $false, $false, $true, 1, $false, "365÷7=52, 1", 2)
$d.Content.Find.Execute("245÷5=49, 0", $true, $false, $false, $false, $false, $true, 1, $false, "437÷3=145, 2", 2)
$d.Content.Find.Execute("127÷7=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "409÷7=58, 3", 2)
$d.Content.Find.Execute("880÷2=440, 0", $true, $false, $false, $false, $false, $true, 1, $false, "320÷3=106, 2", 2)
$d.Content.Find.Execute("410÷8=51, 2", $true, $false, $false, $false, $false, $true, 1, $false, "359÷7=51, 2", 2)
$d.Content.Find.Execute("903÷6=150, 3", $true, $false, $false, $false, $false, $true, 1, $false, "903÷3=301, 0", 2)
$d.Content.Find.Execute("586÷2=293, 0", $true, $false, $false, $false, $false, $true, 1, $false, "479÷2=239, 1", 2)
$d.Content.Find.Execute("803÷2=401, 1", $true, $false, $false, $false, $false, $true, 1, $false, "446÷2=223, 0", 2)
$d.Content.Find.Execute("452÷2=226, 0", $true, $false, $false, $false, $false, $true, 1, $false, "862÷7=123, 1", 2)
$d.Content.Find.Execute("509÷9=56, 5", $true, $false, $false, $false, $false, $true, 1, $false, "436÷5=87, 1", 2)
$d.Content.Find.Execute("864÷8=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "432÷2=216, 0", 2)
$d.Content.Find.Execute("652÷4=163, 0", $true, $false, $false, $false, $false, $true, 1, $false, "579÷8=72, 3", 2)
